# The workbook originally lists the "2022-Q2" detail sheet before the
# "总计" (totals) summary sheet. This commit re-sorts the sheet tabs so the
# "总计" summary sheet comes first, followed by the "2022-Q2" detail sheet,
# without altering any cell data on either sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# Move "总计" so it becomes the first tab (i.e. ahead of "2022-Q2").
$totalSheet.Move($wb.Worksheets.Item(1))
